# Add "OutputAssociateReportPath" / "OutputCohortReportPath" settings rows
# to the Settings sheet, just below the existing OrchestratorQueueName row,
# and push the existing "logF_BusinessProcessName" row down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a single new blank row above the current row 4 (the
# "logF_BusinessProcessName" row). The pre-existing blank row 3 becomes the
# first of the two new settings rows; the freshly inserted row becomes the
# second. Everything from the old row 4 onward (including the old blank
# row 5) shifts down by one row, so the net sheet length grows by exactly
# one row.
$ws.Rows("4:4").Insert()

$ws.Range("A3").Value = "OutputAssociateReportPath"
$ws.Range("C3").Value = "Path to directory for generated associate reports."

$ws.Range("A4").Value = "OutputCohortReportPath"
$ws.Range("C4").Value = "Path to directory for generated cohort reports."

# Match the default (non-wrapped) row height on the two new rows.
$ws.Rows("4:4").RowHeight = 14.25

# The "logF_BusinessProcessName" row (now row 5, wrapped description text)
# recalculates to a taller row height.
$ws.Rows("5:5").RowHeight = 30

# Leave the active selection on the newly added description cell.
$null = $ws.Range("C4").Select()
